# Apply the "appliance type" migration edit to the legacy data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates ---
$ws.Range("A1").Value = "Ime i prezime klijenta"
$ws.Range("D1").Value = "Tip aparata"
$ws.Range("H1").Value = "Opis kvara"

# --- Row 2 updates ---
$ws.Range("D2").Value = "SM"
$ws.Range("F2").Value = "DW50K"
$ws.Range("H2").Value = "ne pere kako treba"

# --- Row 3 updates ---
$ws.Range("D3").Value = "VM"
$ws.Range("F3").Value = "WM-5000"
$ws.Range("H3").Value = "ne centrifugira"

# --- Row 4 updates ---
$ws.Range("D4").Value = "VM KOMB"
$ws.Range("F4").Value = "WKD-300"
$ws.Range("H4").Value = "ne suši"

# --- Row 5 updates ---
$ws.Range("D5").Value = "SM UG"
$ws.Range("F5").Value = "WDI-60"
$ws.Range("H5").Value = "ne radi"

# --- Row 6 updates ---
$ws.Range("H6").Value = "ne hladi"

# --- New row 7 ---
$ws.Range("A7").Value = "Marko Petrović"
$ws.Range("B7").Value = "069/987-654"
$ws.Range("C7").Value = "KO"
$ws.Range("D7").Value = "šporet"
$ws.Range("E7").Value = "Gorenje"
$ws.Range("F7").Value = "G-500"
# G7 looks numeric ("11111"); force text formatting first so it is
# stored as a string (matching the other Serijski broj cells), then
# restore the default "Normal" style so no stray formatting is left
# behind on the cell.
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "11111"
$ws.Range("G7").Style = "Normal"
$ws.Range("H7").Value = "ne radi ploca"
